$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CEMENTO"

$ws.Cells.Item(1,1).Value = "Año"
$ws.Cells.Item(1,2).Value = "Mes"
$ws.Cells.Item(1,3).Value = "Consumo de cemento Zona Oeste"
$ws.Cells.Item(1,4).Value = "Consumo de cemento Zona Oeste. Var interanual"
$ws.Cells.Item(1,5).Value = "Consumo de cemento Zona Oeste. Tendencia"
$ws.Cells.Item(1,6).Value = "Consumo de cemento España"
$ws.Cells.Item(1,7).Value = "Consumo de cemento España. Var interanual"
$ws.Cells.Item(1,8).Value = "Consumo de cemento España. Tendencia"

$ws.Range("D2:E30").NumberFormat = "0.00"
$ws.Range("G2:H30").NumberFormat = "0.00"

$data = @"
2018,1,181909.26748000001,30.185835987371679,5.4173105142018319,1008585,20.728908934862169,9.5901606332262439
2018,2,175855.69600000003,20.093499304617346,5.4077397414733186,992748,7.5284703252235019,9.5490568827034554
2018,3,177703.58094000001,-3.2254551221393979,5.3299807712108418,981175,-12.375608505819603,9.4685831327404202
2018,4,200864.91915999999,21.026213505715052,5.1829766391775323,1093631,19.153982763692625,9.3478077876023402
2018,5,230236.40516000002,12.128926687970566,4.9650762536439279,1251351,9.8620744148478536,9.1842822938017417
2018,6,229352.92300000004,11.533231680373968,4.6757287476629648,1249711,10.230912190918229,8.9762390822244917
2018,7,230310.26400000002,9.8982739014478849,4.3148807439010719,1236481,17.358282831638338,8.721957652653753
2018,8,224144.179,9.027694049653423,3.8829550805061173,1106769,11.387328380410256,8.4198046349496813
2018,9,223875.3835,0.49645727766394199,3.3807623312619102,1130262,7.7957353303030441,8.068746403776526
2018,10,254837.64561999997,12.425715287364625,2.8094703434917832,1270502,12.953491328688372,7.6679554118364175
2018,11,209149.25412,-13.035135187403135,2.1700466655570141,1153494,-4.5470907605902777,7.2165851527291629
2018,12,150874.285111,-6.9003675803792186,1.4641266406066491,986374,10.205075779159479,6.7141561711599032
2019,1,169947.49400000001,-6.5756811874992271,0.6922896963832782,1134694,12.503556963468633,6.1593720898953546
2019,2,207934.88287,18.241767312444622,-0.14546560702474287,1168774,17.73118656496915,5.5511789566750105
2019,3,238027.44600000003,33.946341846857784,-1.0496454297293283,1286754,31.144189364792219,4.8889633876323613
2019,4,217294.18867000003,8.1792627496657175,-2.019479040667429,1202230,9.9301318269142023,4.1729578327625845
2019,5,256713.48300000001,11.499952764464005,-3.0517654318817886,1391711,11.216677015481658,3.4052180216426051
2019,6,217095.68999999997,-5.3442671842468981,-4.1425953494574888,1271376,1.7336008085069299,2.5881994875989407
2019,7,239124.93500000003,3.8273027206464416,-5.2870490034937534,1410564,14.078906186184836,1.7249002263882358
2019,8,209643.60800000001,-6.4693051877113383,-6.4802900535227792,1183473,6.9304434800757875,0.81825888663664204
2019,9,220151.13,-1.6635386355463244,-7.7168492179848096,1235587,9.3186358561112357,-0.12792796594914754
2019,10,232585.53,-8.7318792974492965,-8.991256452482185,1406042,10.668224056317888,-1.1094393083990064
2019,11,173999.405,-16.806107804636341,-10.297621343826798,1104087,-4.2832472470598066,-2.1213981063662759
2019,12,126881.12999999998,-15.902746510678067,-11.63003546652811,886133,-10.162575250361428,-3.1581094322150811
2020,1,154452.50000000003,-9.1175183789411918,-12.983042373322029,1035484,-8.7433263946050666,-4.2140284867220945
2020,2,180033.57499999998,-13.418291094257517,-14.351482332989193,1174390,0.48050350196018332,-5.2840968919013598
2020,3,133937.78378307819,-43.730109265181895,-15.729927175143963,921907,-28.354059905778417,-6.3635708043438566
2020,4,99219.250766992569,-54.338746298606864,-17.11288392445357,598051,-50.25485971902215,-7.4473060611687689
2020,5,173616.34413599968,-32.369604390432549,-18.496804062674826,1143853,-17.80958834125763,-8.5316856167939914
"@

$lines = $data -split "`n"
$r = 2
foreach ($line in $lines) {
    $cols = $line.Trim() -split ","
    $ws.Cells.Item($r, 1).Value = [double]$cols[0]
    $ws.Cells.Item($r, 2).Value = [double]$cols[1]
    $ws.Cells.Item($r, 3).Value = [double]$cols[2]
    $ws.Cells.Item($r, 4).Value = [double]$cols[3]
    $ws.Cells.Item($r, 5).Value = [double]$cols[4]
    $ws.Cells.Item($r, 6).Value = [double]$cols[5]
    $ws.Cells.Item($r, 7).Value = [double]$cols[6]
    $ws.Cells.Item($r, 8).Value = [double]$cols[7]
    $r++
}

$ws.Range("B20").Select()
$ws.Activate()
